# Update Name of Algo
# Applies updated KNN imputation results to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.030000000000001
$ws.Range("E3").Value = 16.751
$ws.Range("B12").Value = 5.220999999999999
$ws.Range("D14").Value = -7.547
$ws.Range("D26").Value = -7.699000000000001
$ws.Range("E30").Value = 15.899
$ws.Range("D31").Value = -8.430000000000001
$ws.Range("B32").Value = 6.267
$ws.Range("D35").Value = -7.672
$ws.Range("B36").Value = 8.548
$ws.Range("D37").Value = -7.741
$ws.Range("B38").Value = 5.445
$ws.Range("E44").Value = 16.747
$ws.Range("D45").Value = -7.498
$ws.Range("B46").Value = 6.377000000000001
$ws.Range("B54").Value = 5.155
$ws.Range("B55").Value = 4.782999999999999
$ws.Range("D57").Value = -8.101000000000001
$ws.Range("E58").Value = 16.425
$ws.Range("B67").Value = 5.329
$ws.Range("B69").Value = 5.147
$ws.Range("B72").Value = 5.380999999999999
$ws.Range("E84").Value = 16.3
$ws.Range("E89").Value = 17.199
$ws.Range("B91").Value = 5.276
$ws.Range("E91").Value = 16.916
$ws.Range("E92").Value = 16.931
$ws.Range("B99").Value = 5.217000000000001
$ws.Range("D100").Value = -8.280000000000001
$ws.Range("D102").Value = -7.808
$ws.Range("E102").Value = 16.692
